$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 4
$ws.Range("C2").Value = 11
$ws.Range("B3").Value = 5
$ws.Range("C3").Value = 9.5
$ws.Range("C5").Value = 20
